$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 44830
$ws.Range("A1").NumberFormat = "d-mmm-yy"

$ws.Range("B1").Value = "https://github.com/gandharvas/crs/blob/main/files/26_09_2022.xlsx?raw=true"

$ws.Range("B1").Select()

$ws.Columns.Item(1).AutoFit() | Out-Null
